# Employee training report update:
#  - Title and header row font becomes bold + white (no longer a distinct
#    larger 14pt font for the title; it now shares the same bold/white
#    font used by the header row).
#  - "PERIOD TO EXPIRE" (column H) values drop by 8 days.
#  - "LAST UPDATE" (column I) moves from 08-Sep-2025 to 16-Sep-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Formatting: title (A1) and header row (A2:K2) share one bold white font ---
$titleRange  = $ws.Range("A1")
$headerRange = $ws.Range("A2:K2")

$titleRange.Font.Bold = $true
$titleRange.Font.Size = 11
$titleRange.Font.Color = 16777215   # white (RGB 255,255,255)

$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215  # white (RGB 255,255,255)

# --- Data updates: PERIOD TO EXPIRE (H) and LAST UPDATE (I) for rows 3-17 ---
$periodToExpire = @{
    3  = 672
    4  = 674
    5  = 672
    6  = 674
    7  = 672
    8  = 673
    9  = 674
    10 = 673
    11 = 674
    12 = 675
    13 = 675
    14 = 675
    15 = 308
    16 = 314
    17 = 314
}

# Mark the LAST UPDATE column as text first so the "dd-mmm-yyyy" looking
# string is preserved literally instead of being auto-converted into a
# date serial number.
$lastUpdateRange = $ws.Range("I3:I17")
$lastUpdateRange.NumberFormat = "@"

foreach ($row in 3..17) {
    $ws.Range("H$row").Value = $periodToExpire[$row]
    $ws.Range("I$row").Value = "16-Sep-2025"
}
